# ndimas add date fix bug
#
# Add a "date" column (G) to the accounts sheet, one value per existing
# account row. The dates are entered with a leading apostrophe so Excel
# stores them as literal text ("d/m/yyyy") instead of re-parsing them into
# serial date numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new values in the same order the rows already appear on the
# sheet (row1, row4, row2, row3) so the shared-string table grows in that
# sequence.
$ws.Range("G1").Value = "'1/1/2000"
$ws.Range("G4").Value = "'12/12/2000"
$ws.Range("G2").Value = "'15/3/2000"
$ws.Range("G3").Value = "'12/1/2000"

# Size the new column to comfortably fit the date text.
$ws.Columns("G:G").ColumnWidth = 9.6

# Leave the selection where the user ended up after typing the last value.
[void]$ws.Range("G4").Select()
